$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -77.0651
$ws.Range("B2").Value = -76.9564

$ws.Range("A3").Value = 37.0547
$ws.Range("B3").Value = 37.1403

$ws.Range("A4").Value = -76.2058
$ws.Range("B4").Value = -76.3155

$ws.Range("A5").Value = 37.7271
$ws.Range("B5").Value = 37.6418
